$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style cleanup -----------------------------------------------------
# The old style index 2 (numFmt-flavoured "center/center") is dropped; the
# A-column counter cells move to the plain "center/center" style (index 1),
# which also happens to be column A's own default style. Re-asserting the
# alignment on the whole A2:A16 block makes the engine resolve the cells to
# the existing equivalent xf instead of keeping the now-redundant one.
# (-4108 = xlCenter)
$ws.Range("A2:A16").HorizontalAlignment = -4108
$ws.Range("A2:A16").VerticalAlignment = -4108

# D2/D3 used the old "left/center" style (old index 3); re-apply the same
# alignment so they resolve to the shifted index after the prune above.
# (-4131 = xlLeft, -4108 = xlCenter)
$ws.Range("D2:D3").HorizontalAlignment = -4131
$ws.Range("D2:D3").VerticalAlignment = -4108

# --- Refresh the "No" counter formulas ---------------------------------
# Re-assigning the same formula across A2:A7 as one Range.Formula write
# regroups it into its own shared-formula block (separate from A8:A16,
# which already is one), matching the split the author ended up with.
$ws.Range("A2:A7").Formula = "=ROW()-1"
$ws.Range("A8:A16").Formula = "=ROW()-1"

# --- New TODO rows -------------------------------------------------------
$ws.Range("B14").Value = "TODO"
$ws.Range("C14").Value = "首页的scam case列表添加下拉刷新功能"
$ws.Range("D14").Value = "在添加post之后下拉列表刷新，主页显示新添加的案例"
$ws.Range("E14").Value = "未测试"
$ws.Range("F14").Value = "Yijing"

$ws.Range("B15").Value = "TODO"
$ws.Range("C15").Value = "Profile添加修改照片和名字的功能"
$ws.Range("D15").Value = "点击头像，能够上传图片，点击名字可以修改"
$ws.Range("E15").Value = "未测试"
$ws.Range("F15").Value = "Yijing"

$ws.Range("B16").Value = "TODO"
$ws.Range("C16").Value = "重构firebaseauth和firestore的代码，形成单例模式和state模式"
$ws.Range("D16").Value = "相关所有功能正常运行"
$ws.Range("E16").Value = "未测试"
$ws.Range("F16").Value = "Zhaoyun"

$ws.Range("B17").Value = "TODO"
$ws.Range("C17").Value = "用户读取数据和添加数据加入权限"
$ws.Range("D17").Value = "相关所有功能正常运行"
$ws.Range("E17").Value = "未测试"
$ws.Range("F17").Value = "Zhaoyun"
